$d = $word.ActiveDocument

# 1. "Your model includes " -> "model includes "
$d.Content.Find.Execute("Your model includes ", $true, $false, $false, $false, $false, $true, 1, $false, "model includes ", 2) | Out-Null

# 2. Remove the two inline images (Picture 1 / rId13 and Picture 2 / rId14).
#    Deleting the shape itself (not its Range) removes only the drawing run
#    and leaves the paragraph (and its pPr, if any) intact.
while ($d.InlineShapes.Count -gt 0) {
    $d.InlineShapes.Item(1).Delete()
}

# 3. Remove the now-empty "Apple Color Emoji" paragraph and the
#    "Suggested Wording for Your Report" Heading2 paragraph that followed
#    the second (now-removed) image.
$regionIdx = -1
$idx = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "Region Defined by Spatial Grid") {
        $regionIdx = $idx
    }
    $idx = $idx + 1
}
# Layout right after "Region Defined...": [image-2 para] [blank para] [Heading2 para] [next content]
# Keep the image-2 paragraph (now empty, still carrying its pPr) and delete the
# following blank paragraph and the Heading2 paragraph entirely (each paragraph's
# own Range includes its paragraph mark, so deleting it removes the whole
# paragraph rather than just merging text).
$blankPara = $d.Paragraphs.Item($regionIdx + 2)
$d.Range($blankPara.Range.Start, $blankPara.Range.End).Delete()
$headingPara = $d.Paragraphs.Item($regionIdx + 2)
$d.Range($headingPara.Range.Start, $headingPara.Range.End).Delete()

# 4. Remove the stray lastRenderedPageBreak on the "First, we incorporated..." run
#    by deleting the paragraph's content (keeping the paragraph mark / pPr) and
#    retyping the text.
$firstIdx = -1
$idx = 1
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "First, we incorporated temporal dependence") {
        $firstIdx = $idx
    }
    $idx = $idx + 1
}
$p = $d.Paragraphs.Item($firstIdx)
$start = $p.Range.Start
$end = $p.Range.End - 1
$d.Range($start, $end).Delete()
$d.Range($start, $start).InsertAfter("First, we incorporated temporal dependence by including lagged average magnitude as a covariate.")

# 5. Merge ", RMS, " + "Clo" + "), rather..." runs (dropping the proofErr tags
#    around "Clo") into a single run of text.
$d.Content.Find.Execute(", RMS, Clo), rather than treating it as a global or region-specific scalar.", $true, $false, $false, $false, $false, $true, 1, $false, ", RMS, Clo), rather than treating it as a global or region-specific scalar.", 2) | Out-Null
